$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that parse as plain numbers need NumberFormat forced to
# Text before assignment so Excel keeps them as strings (matching the
# original inline-string cell type), then the temporary format is cleared
# so the cell keeps the workbook default style (no explicit s= attribute).
$textForcedCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D14', 'D16', 'D18', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D32', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.093.92'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.353.09'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '542.11'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').Value = '134.17'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +5.51%  '
$ws.Range('D9').Value = '0.105'
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('D10').Value = '5.55'
$ws.Range('E10').Value = '  +3.00%  '
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.771.20'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '23.76'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = '58.042.42'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '0.0000136'
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').Value = '2.335.25'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '10.80'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').Value = '4.31'
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').Value = '329.87'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').Value = '63.27'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '0.165'
$ws.Range('E24').Value = '  -2.83%  '
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = '8.29'
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('E27').Value = '  -5.27%  '
$ws.Range('D28').Value = '1.76'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '170.30'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '0.0₃0737'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '18.37'
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').Value = '4.20'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '1.24'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').Value = '39.06'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').Value = '0.379'
$ws.Range('D41').Value = '290.66'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.65'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '140.92'
$ws.Range('E43').Value = '  -6.17%  '
$ws.Range('D44').Value = '0.0950'
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('D45').Value = '0.0512'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.567'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '18.93'
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('D48').Value = '0.0223'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('D49').Value = '0.382'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').Value = '4.70'
$ws.Range('E51').Value = '  +0.70%  '

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).ClearFormats()
}
